# Updated cryptos list on Mon Jun 19 02:06:59 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be stored as text so values such as "1.000",
# "0.9990", "26.435.33" etc. are preserved exactly (and are not silently
# re-interpreted by Excel as numbers/dates).
$ws.Range("D2:D51").NumberFormat = "@"

# --- Rows 2-33: only the Price (D) and Volume(1h) (E) columns changed ---

$ws.Range("D2").Value  = "26.435.33"
$ws.Range("E2").Value  = "  -0.17%  "

$ws.Range("D3").Value  = "1.724.38"
$ws.Range("E3").Value  = "  -0.16%  "

$ws.Range("D4").Value  = "0.9990"
$ws.Range("E4").Value  = "  +0.00%  "

$ws.Range("D5").Value  = "243.18"
$ws.Range("E5").Value  = "  -0.69%  "

$ws.Range("D6").Value  = "0.9995"
$ws.Range("E6").Value  = "  +0.00%  "

$ws.Range("D7").Value  = "0.4900"
$ws.Range("E7").Value  = "  +2.02%  "

$ws.Range("D8").Value  = "0.2608"
$ws.Range("E8").Value  = "  -2.66%  "

$ws.Range("D9").Value  = "0.06190"
$ws.Range("E9").Value  = "  +0.02%  "

$ws.Range("D10").Value = "1.728.00"
$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("D11").Value = "0.06983"
$ws.Range("E11").Value = "  -2.57%  "

$ws.Range("D12").Value = "15.52"
$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("D13").Value = "4.557"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("D14").Value = "0.5988"
$ws.Range("E14").Value = "  -1.58%  "

$ws.Range("D15").Value = "77.32"
$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("D16").Value = "0.9995"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Value = "26.435.74"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").Value = "0.9992"
$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").Value = "0.000007182"
$ws.Range("E19").Value = "  +3.31%  "

# Row 20: only Volume(1h) changed, Price stays "11.32"
$ws.Range("E20").Value = "  -1.45%  "

$ws.Range("D21").Value = "1.948.70"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").Value = "4.459"
$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("D23").Value = "8.571"
$ws.Range("E23").Value = "  -2.28%  "

$ws.Range("D24").Value = "5.160"
$ws.Range("E24").Value = "  -1.79%  "

$ws.Range("D25").Value = "137.12"
$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("D26").Value = "15.24"
$ws.Range("E26").Value = "  -0.45%  "

$ws.Range("D27").Value = "1.399"
$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("D28").Value = "107.04"
$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("D29").Value = "1.707"
$ws.Range("E29").Value = "  -4.14%  "

$ws.Range("D30").Value = "3.954"
$ws.Range("E30").Value = "  -0.36%  "

$ws.Range("D31").Value = "0.07972"
$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("D32").Value = "3.666"
$ws.Range("E32").Value = "  -0.45%  "

$ws.Range("D33").Value = "0.04500"
$ws.Range("E33").Value = "  -0.35%  "

# --- Rows 34-51: a new coin ("Frax") was inserted at row 34, which shifts
#     all the following coins down by one row and drops the last one
#     (Aave) off the bottom of the list. ---

$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").Value = "0.9988"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.604"
$ws.Range("E35").Value = "  -0.42%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  +0.86%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.6234"
$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "0.9040"
$ws.Range("E38").Value = "  -0.69%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "1.964"
$ws.Range("E39").Value = "  -5.51%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.391"
$ws.Range("E40").Value = "  +0.67%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "0.9997"
$ws.Range("E41").Value = "  -0.61%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.01481"
$ws.Range("E42").Value = "  -1.36%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "99.92"
$ws.Range("E43").Value = "  -3.80%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.447"
$ws.Range("E44").Value = "  -2.81%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.3844"
$ws.Range("E45").Value = "  -0.49%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "6.693"
$ws.Range("E46").Value = "  -3.06%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.1152"
$ws.Range("E47").Value = "  -2.27%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05362"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.713"
$ws.Range("E49").Value = "  -0.92%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "30.06"
$ws.Range("E50").Value = "  -1.46%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.238"
$ws.Range("E51").Value = "  -1.18%  "
